$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.651.68'
$ws.Range("E2").Value = '  -0.79%  '

$ws.Range("D3").Value = '2.461.35'
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''574.09'
$ws.Range("E5").Value = '  -0.73%  '

$ws.Range("D6").Value = '''146.96'
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  -1.61%  '

$ws.Range("E9").Value = '  -0.86%  '

$ws.Range("E10").Value = '  -0.66%  '

$ws.Range("E11").Value = '  -0.74%  '

$ws.Range("E12").Value = '  -0.64%  '

$ws.Range("D13").Value = '''28.98'
$ws.Range("E13").Value = '  +1.86%  '

$ws.Range("E14").Value = '  -2.53%  '

$ws.Range("D15").Value = '2.909.05'
$ws.Range("E15").Value = '  +0.35%  '

$ws.Range("D16").Value = '62.606.01'
$ws.Range("E16").Value = '  -0.82%  '

$ws.Range("D17").Value = '2.461.07'
$ws.Range("E17").Value = '  -0.02%  '

$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("D19").Value = '''10.93'
$ws.Range("E19").Value = '  -1.42%  '

$ws.Range("D20").Value = '''325.20'
$ws.Range("E20").Value = '  -1.77%  '

$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("D22").Value = '''2.18'
$ws.Range("E22").Value = '  +1.54%  '

$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").Value = '''10.02'
$ws.Range("E24").Value = '  +16.04%  '

$ws.Range("D25").Value = '''65.34'
$ws.Range("E25").Value = '  -1.65%  '

$ws.Range("D26").Value = '''638.24'
$ws.Range("E26").Value = '  -2.01%  '

$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("D28").Value = '0.0₃0972'
$ws.Range("E28").Value = '  -4.43%  '

$ws.Range("E29").Value = '  -20.40%  '

$ws.Range("E30").Value = '  -1.09%  '

$ws.Range("D31").Value = '''7.92'
$ws.Range("E31").Value = '  -3.70%  '

$ws.Range("E32").Value = '  -2.70%  '

$ws.Range("E33").Value = '  -2.87%  '

$ws.Range("E34").Value = '  -0.10%  '

$ws.Range("E35").Value = '  +2.10%  '

$ws.Range("E36").Value = '  -1.56%  '

$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = '''0.368'
$ws.Range("E37").Value = '  -1.90%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = '''151.16'
$ws.Range("E38").Value = '  -0.85%  '

$ws.Range("D39").Value = '''18.60'
$ws.Range("E39").Value = '  -1.47%  '

$ws.Range("D40").Value = '''5.31'

$ws.Range("D41").Value = '''2.73'
$ws.Range("E41").Value = '  -1.00%  '

$ws.Range("E42").Value = '  -2.88%  '

$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("E44").Value = '  -23.71%  '

$ws.Range("D45").Value = '''153.11'
$ws.Range("E45").Value = '  +3.97%  '

$ws.Range("D46").Value = '''15.29'
$ws.Range("E46").Value = '  +1.93%  '

$ws.Range("E47").Value = '  -2.06%  '

$ws.Range("E48").Value = '  -0.33%  '

$ws.Range("D49").Value = '''20.26'
$ws.Range("E49").Value = '  -2.46%  '

$ws.Range("D50").Value = '''0.0507'
$ws.Range("E50").Value = '  -2.02%  '

$ws.Range("D51").Value = '''0.0908'
$ws.Range("E51").Value = '  -1.73%  '
